$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 41.92137399999999
$ws.Range("H2").Value = 125.764122
$ws.Range("I2").Value = 0.2188311536698969
$ws.Range("J2").Value = 0.2273746866916212
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.107177
$ws.Range("N2").Value = 0.321531
$ws.Range("O2").Value = 0.003526763356587491
$ws.Range("P2").Value = 0.003549676734010809
$ws.Range("Q2").Value = 4.493007101198
$ws.Range("R2").Value = 40.43706391078199
$ws.Range("S2").Value = 0.0007717656940427586
$ws.Range("T2").Value = 0.0008071066352522447
$ws.Range("G3").Value = 41.92137399999999
$ws.Range("H3").Value = 125.764122
$ws.Range("I3").Value = 0.2188311536698969
$ws.Range("J3").Value = 0.2273746866916212
$ws.Range("O3").Value = 0.9757678722356318
$ws.Range("P3").Value = 0.9821074349659524
$ws.Range("Q3").Value = 1243.103530291201
$ws.Range("R3").Value = 11187.93177262081
$ws.Range("S3").Value = 0.2135284091953439
$ws.Range("T3").Value = 0.2233063703228952
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 41.92137399999999
$ws.Range("H4").Value = 125.764122
$ws.Range("I4").Value = 0.2188311536698969
$ws.Range("J4").Value = 0.2273746866916212
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.04072766666666667
$ws.Range("N4").Value = 0.122183
$ws.Range("O4").Value = 0.001340183457265176
$ws.Range("P4").Value = 0.001348890627627329
$ws.Range("Q4").Value = 1.707359746480666
$ws.Range("R4").Value = 15.366237718326
$ws.Range("S4").Value = 0.0002932738920826495
$ws.Range("T4").Value = 0.0003067035838380281
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 41.92137399999999
$ws.Range("H5").Value = 125.764122
$ws.Range("I5").Value = 0.2188311536698969
$ws.Range("J5").Value = 0.2273746866916212
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5885005
$ws.Range("N5").Value = 1.177001
$ws.Range("O5").Value = 0.01936518095051565
$ws.Range("P5").Value = 0.01299399767240936
$ws.Range("Q5").Value = 24.670749559687
$ws.Range("R5").Value = 148.024497358122
$ws.Range("S5").Value = 0.004237704888427651
$ws.Range("T5").Value = 0.002954506149635734
$ws.Range("I6").Value = 0.1210191186482915
$ws.Range("J6").Value = 0.1257439067741692
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.107177
$ws.Range("N6").Value = 0.321531
$ws.Range("O6").Value = 0.003526763356587491
$ws.Range("P6").Value = 0.003549676734010809
$ws.Range("Q6").Value = 2.484745660518334
$ws.Range("R6").Value = 22.362710944665
$ws.Range("S6").Value = 0.0004268057930953082
$ws.Range("T6").Value = 0.0004463502203198924
$ws.Range("I7").Value = 0.1210191186482915
$ws.Range("J7").Value = 0.1257439067741692
$ws.Range("O7").Value = 0.9757678722356318
$ws.Range("P7").Value = 0.9821074349659524
$ws.Range("S7").Value = 0.1180865679032748
$ws.Range("T7").Value = 0.1234940257445771
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("I8").Value = 0.1210191186482915
$ws.Range("J8").Value = 0.1257439067741692
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.04072766666666667
$ws.Range("N8").Value = 0.122183
$ws.Range("O8").Value = 0.001340183457265176
$ws.Range("P8").Value = 0.001348890627627329
$ws.Range("Q8").Value = 0.9442127789827778
$ws.Range("R8").Value = 8.497915010845
$ws.Range("S8").Value = 0.0001621878208252518
$ws.Range("T8").Value = 0.0001696147773289213
$ws.Range("D9").Value = "MuSCs"
$ws.Range("I9").Value = 0.1210191186482915
$ws.Range("J9").Value = 0.1257439067741692
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.5885005
$ws.Range("N9").Value = 1.177001
$ws.Range("O9").Value = 0.01936518095051565
$ws.Range("P9").Value = 0.01299399767240936
$ws.Range("Q9").Value = 13.64354351761917
$ws.Range("R9").Value = 81.86126110571499
$ws.Range("S9").Value = 0.002343557131096088
$ws.Range("T9").Value = 0.001633916031943214
$ws.Range("G10").Value = 50.59004100000001
$ws.Range("H10").Value = 151.770123
$ws.Range("I10").Value = 0.2640819224159348
$ws.Range("J10").Value = 0.274392120880658
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.107177
$ws.Range("N10").Value = 0.321531
$ws.Range("O10").Value = 0.003526763356587491
$ws.Range("P10").Value = 0.003549676734010809
$ws.Range("Q10").Value = 5.422088824257001
$ws.Range("R10").Value = 48.798799418313
$ws.Range("S10").Value = 0.0009313544471136996
$ws.Range("T10").Value = 0.0009740033274859529
$ws.Range("G11").Value = 50.59004100000001
$ws.Range("H11").Value = 151.770123
$ws.Range("I11").Value = 0.2640819224159348
$ws.Range("J11").Value = 0.274392120880658
$ws.Range("O11").Value = 0.9757678722356318
$ws.Range("P11").Value = 0.9821074349659524
$ws.Range("Q11").Value = 1500.157379495162
$ws.Range("R11").Value = 13501.41641545646
$ws.Range("S11").Value = 0.2576826555316919
$ws.Range("T11").Value = 0.2694825420129706
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("G12").Value = 50.59004100000001
$ws.Range("H12").Value = 151.770123
$ws.Range("I12").Value = 0.2640819224159348
$ws.Range("J12").Value = 0.274392120880658
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.04072766666666667
$ws.Range("N12").Value = 0.122183
$ws.Range("O12").Value = 0.001340183457265176
$ws.Range("P12").Value = 0.001348890627627329
$ws.Range("Q12").Value = 2.060414326501
$ws.Range("R12").Value = 18.543728938509
$ws.Range("S12").Value = 0.0003539182237846216
$ws.Range("T12").Value = 0.0003701249601507045
$ws.Range("D13").Value = "MuSCs"
$ws.Range("G13").Value = 50.59004100000001
$ws.Range("H13").Value = 151.770123
$ws.Range("I13").Value = 0.2640819224159348
$ws.Range("J13").Value = 0.274392120880658
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.5885005
$ws.Range("N13").Value = 1.177001
$ws.Range("O13").Value = 0.01936518095051565
$ws.Range("P13").Value = 0.01299399767240936
$ws.Range("Q13").Value = 29.7722644235205
$ws.Range("R13").Value = 178.633586541123
$ws.Range("S13").Value = 0.005113994213344614
$ws.Range("T13").Value = 0.003565450580050738
$ws.Range("G14").Value = 21.5944925
$ws.Range("H14").Value = 43.188985
$ws.Range("I14").Value = 0.1127240654538407
$ws.Range("J14").Value = 0.07808333391699843
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.107177
$ws.Range("N14").Value = 0.321531
$ws.Range("O14").Value = 0.003526763356587491
$ws.Range("P14").Value = 0.003549676734010809
$ws.Range("Q14").Value = 2.3144329226725
$ws.Range("R14").Value = 13.886597536035
$ws.Range("S14").Value = 0.0003975511034481754
$ws.Range("T14").Value = 0.0002771705937191664
$ws.Range("G15").Value = 21.5944925
$ws.Range("H15").Value = 43.188985
$ws.Range("I15").Value = 0.1127240654538407
$ws.Range("J15").Value = 0.07808333391699843
$ws.Range("O15").Value = 0.9757678722356318
$ws.Range("P15").Value = 0.9821074349659524
$ws.Range("Q15").Value = 640.3461361165516
$ws.Range("R15").Value = 3842.07681669931
$ws.Range("S15").Value = 0.1099925214976443
$ws.Range("T15").Value = 0.07668622278681328
$ws.Range("D16").Value = "Inflammatory-Mac"
$ws.Range("G16").Value = 21.5944925
$ws.Range("H16").Value = 43.188985
$ws.Range("I16").Value = 0.1127240654538407
$ws.Range("J16").Value = 0.07808333391699843
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.04072766666666667
$ws.Range("N16").Value = 0.122183
$ws.Range("O16").Value = 0.001340183457265176
$ws.Range("P16").Value = 0.001348890627627329
$ws.Range("Q16").Value = 0.8794932923758334
$ws.Range("R16").Value = 5.276959754255
$ws.Range("S16").Value = 0.0001510709277569143
$ws.Range("T16").Value = 0.0001053258772945343
$ws.Range("D17").Value = "MuSCs"
$ws.Range("G17").Value = 21.5944925
$ws.Range("H17").Value = 43.188985
$ws.Range("I17").Value = 0.1127240654538407
$ws.Range("J17").Value = 0.07808333391699843
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.5885005
$ws.Range("N17").Value = 1.177001
$ws.Range("O17").Value = 0.01936518095051565
$ws.Range("P17").Value = 0.01299399767240936
$ws.Range("Q17").Value = 12.70836963349625
$ws.Range("R17").Value = 50.833478533985
$ws.Range("S17").Value = 0.002182921924991396
$ws.Range("T17").Value = 0.001014614659171441
$ws.Range("G18").Value = 54.28001766666667
$ws.Range("H18").Value = 162.840053
$ws.Range("I18").Value = 0.2833437398120361
$ws.Range("J18").Value = 0.2944059517365532
$ws.Range("K18").Value = 2
$ws.Range("L18").Value = 0.6666666666666666
$ws.Range("M18").Value = 0.107177
$ws.Range("N18").Value = 0.321531
$ws.Range("O18").Value = 0.003526763356587491
$ws.Range("P18").Value = 0.003549676734010809
$ws.Range("Q18").Value = 5.817569453460335
$ws.Range("R18").Value = 52.35812508114301
$ws.Range("S18").Value = 0.000999286318887549
$ws.Range("T18").Value = 0.001045045957233552
$ws.Range("G19").Value = 54.28001766666667
$ws.Range("H19").Value = 162.840053
$ws.Range("I19").Value = 0.2833437398120361
$ws.Range("J19").Value = 0.2944059517365532
$ws.Range("O19").Value = 0.9757678722356318
$ws.Range("P19").Value = 0.9821074349659524
$ws.Range("Q19").Value = 1609.577052166804
$ws.Range("R19").Value = 14486.19346950124
$ws.Range("S19").Value = 0.2764777181076769
$ws.Range("T19").Value = 0.2891382740986962
$ws.Range("D20").Value = "Inflammatory-Mac"
$ws.Range("G20").Value = 54.28001766666667
$ws.Range("H20").Value = 162.840053
$ws.Range("I20").Value = 0.2833437398120361
$ws.Range("J20").Value = 0.2944059517365532
$ws.Range("L20").Value = 0.6666666666666666
$ws.Range("M20").Value = 0.04072766666666667
$ws.Range("N20").Value = 0.122183
$ws.Range("O20").Value = 0.001340183457265176
$ws.Range("P20").Value = 0.001348890627627329
$ws.Range("Q20").Value = 2.210698466188778
$ws.Range("R20").Value = 19.896286195699
$ws.Range("S20").Value = 0.0003797325928157391
$ws.Range("T20").Value = 0.0003971214290151403
$ws.Range("D21").Value = "MuSCs"
$ws.Range("G21").Value = 54.28001766666667
$ws.Range("H21").Value = 162.840053
$ws.Range("I21").Value = 0.2833437398120361
$ws.Range("J21").Value = 0.2944059517365532
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.5885005
$ws.Range("N21").Value = 1.177001
$ws.Range("O21").Value = 0.01936518095051565
$ws.Range("P21").Value = 0.01299399767240936
$ws.Range("Q21").Value = 31.94381753684217
$ws.Range("R21").Value = 191.662905221053
$ws.Range("S21").Value = 0.005487002792655906
$ws.Range("T21").Value = 0.003825510251608236

Write-Host "Applied all changes"